$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 883.3125
$ws.Range("I6").Value = 1099
$ws.Range("K6").Value = 3297
$ws.Range("M6").Value = -3185
$ws.Range("H38").Value = 3663.5
$ws.Range("I38").Value = 3599
$ws.Range("K38").Value = 10797
$ws.Range("M38").Value = -10425
$ws.Range("H116").Value = 21667.334
$ws.Range("I116").Value = 4999.6665
$ws.Range("K116").Value = 4999.6665
$ws.Range("M116").Value = -1557.6665
$ws.Range("H126").Value = 80912.17999999999
$ws.Range("J126").Value = 80912.17999999999
$ws.Range("L126").Value = 80912.17999999999
$ws.Range("N126").Value = -90792.17999999999
$ws.Range("H127").Value = 200
$ws.Range("I127").Value = 200
$ws.Range("K127").Value = 600
$ws.Range("M127").Value = 4360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3873.232
$ws.Range("I32").Value = 3542.6365
$ws.Range("K32").Value = 3542.6365
$ws.Range("M32").Value = -3255.6365
$ws.Range("H45").Value = 2477.7856
$ws.Range("J45").Value = 2933
$ws.Range("L45").Value = 2933
$ws.Range("N45").Value = -3687
$ws.Range("H61").Value = 2040.8214
$ws.Range("I61").Value = 1698.12
$ws.Range("J61").Value = 4896.6665
$ws.Range("K61").Value = 1698.12
$ws.Range("L61").Value = 4896.6665
$ws.Range("M61").Value = -1486.12
$ws.Range("N61").Value = -5320.6665
$ws.Range("H97").Value = 1394.25
$ws.Range("I97").Value = 1394.25
$ws.Range("K97").Value = 1394.25
$ws.Range("M97").Value = -898.25
$ws.Range("H132").Value = 2273.62
$ws.Range("I132").Value = 2179.422
$ws.Range("J132").Value = 3121.4
$ws.Range("K132").Value = 6538.266
$ws.Range("L132").Value = 9364.200000000001
$ws.Range("M132").Value = -4008.266
$ws.Range("N132").Value = -14424.2
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 2040.8214
$ws.Range("I136").Value = 1698.12
$ws.Range("J136").Value = 4896.6665
$ws.Range("K136").Value = 5094.36
$ws.Range("L136").Value = 14689.9995
$ws.Range("M136").Value = -2544.36
$ws.Range("N136").Value = -19789.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 769.3333
$ws.Range("J64").Value = 811.5
$ws.Range("L64").Value = 811.5
$ws.Range("N64").Value = -1261.5
$ws.Range("H67").Value = 769.3333
$ws.Range("J67").Value = 811.5
$ws.Range("L67").Value = 811.5
$ws.Range("N67").Value = -2371.5
$ws.Range("H86").Value = 1891249.9
$ws.Range("I86").Value = 2835541.5
$ws.Range("K86").Value = 2835541.5
$ws.Range("M86").Value = -2834418.5
$ws.Range("H89").Value = 1891249.9
$ws.Range("I89").Value = 2835541.5
$ws.Range("K89").Value = 14177707.5
$ws.Range("M89").Value = -14172091.5
$ws.Range("H134").Value = 51966.76
$ws.Range("I134").Value = 4455.625
$ws.Range("K134").Value = 13366.875
$ws.Range("M134").Value = -10831.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 125057.5
$ws.Range("J57").Value = 130060
$ws.Range("L57").Value = 130060
$ws.Range("N57").Value = -131180
$ws.Range("H58").Value = 7187.8276
$ws.Range("I58").Value = 6977.095
$ws.Range("J58").Value = 7741
$ws.Range("K58").Value = 6977.095
$ws.Range("L58").Value = 7741
$ws.Range("M58").Value = -6774.095
$ws.Range("N58").Value = -8147
$ws.Range("H99").Value = 5083
$ws.Range("I99").Value = 4282.8
$ws.Range("J99").Value = 6416.6665
$ws.Range("K99").Value = 4282.8
$ws.Range("L99").Value = 6416.6665
$ws.Range("M99").Value = -2784.8
$ws.Range("N99").Value = -9412.666499999999
$ws.Range("H126").Value = 5083
$ws.Range("I126").Value = 4282.8
$ws.Range("J126").Value = 6416.6665
$ws.Range("K126").Value = 12848.4
$ws.Range("L126").Value = 19249.9995
$ws.Range("M126").Value = -10378.4
$ws.Range("N126").Value = -24189.9995
$ws.Range("H132").Value = 1953.1765
$ws.Range("I132").Value = 1950
$ws.Range("K132").Value = 5850
$ws.Range("M132").Value = -3320
$ws.Range("H134").Value = 480023.1
$ws.Range("I134").Value = 4216.1113
$ws.Range("K134").Value = 12648.3339
$ws.Range("M134").Value = -10113.3339
$ws.Range("H136").Value = 7187.8276
$ws.Range("I136").Value = 6977.095
$ws.Range("J136").Value = 7741
$ws.Range("K136").Value = 20931.285
$ws.Range("L136").Value = 23223
$ws.Range("M136").Value = -18381.285
$ws.Range("N136").Value = -28323
$ws.Range("H141").Value = 406584
$ws.Range("J141").Value = 454298.66
$ws.Range("L141").Value = 454298.66
$ws.Range("N141").Value = -464658.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 95.5
$ws.Range("I6").Value = 47.894737
$ws.Range("K6").Value = 143.684211
$ws.Range("M6").Value = -30.684211
$ws.Range("H68").Value = 5000400
$ws.Range("I68").Value = 10000002
$ws.Range("J68").Value = 3333866
$ws.Range("K68").Value = 30000006
$ws.Range("L68").Value = 10001598
$ws.Range("M68").Value = -29999195
$ws.Range("N68").Value = -10003220
$ws.Range("H71").Value = 5000400
$ws.Range("I71").Value = 10000002
$ws.Range("J71").Value = 3333866
$ws.Range("K71").Value = 90000018
$ws.Range("L71").Value = 30004794
$ws.Range("M71").Value = -89995962
$ws.Range("N71").Value = -30012906
$ws.Range("H114").Value = 142857470
$ws.Range("I114").Value = 333333540
$ws.Range("J114").Value = 437.5
$ws.Range("K114").Value = 1000000620
$ws.Range("L114").Value = 1312.5
$ws.Range("M114").Value = -999997366
$ws.Range("N114").Value = -7820.5
$ws.Range("H117").Value = 908.4545000000001
$ws.Range("I117").Value = 658
$ws.Range("J117").Value = 1002.375
$ws.Range("K117").Value = 1974
$ws.Range("L117").Value = 3007.125
$ws.Range("M117").Value = 1468
$ws.Range("N117").Value = -9891.125
$ws.Range("H122").Value = 112635.445
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 144531.28
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 1300781.52
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -1305681.52
$ws.Range("H129").Value = 55574070
$ws.Range("I129").Value = 111112800
$ws.Range("J129").Value = 35344.332
$ws.Range("K129").Value = 333338400
$ws.Range("L129").Value = 106032.996
$ws.Range("M129").Value = -333333400
$ws.Range("N129").Value = -116032.996
$ws.Range("H132").Value = 787918.1
$ws.Range("J132").Value = 2004440.8
$ws.Range("L132").Value = 18039967.2
$ws.Range("N132").Value = -18045027.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H42").Value = 79289
$ws.Range("J42").Value = 79289
$ws.Range("L42").Value = 79289
$ws.Range("N42").Value = -80259
$ws.Range("H97").Value = 2361.8
$ws.Range("I97").Value = 1702.25
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 1702.25
$ws.Range("L97").Value = 5000
$ws.Range("M97").Value = -1206.25
$ws.Range("N97").Value = -5992
$ws.Range("H115").Value = 79289
$ws.Range("J115").Value = 79289
$ws.Range("L115").Value = 79289
$ws.Range("N115").Value = -81639

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 513500
$ws.Range("J20").Value = 513500
$ws.Range("L20").Value = 513500
$ws.Range("N20").Value = -513952
$ws.Range("H36").Value = 35265
$ws.Range("J36").Value = 35265
$ws.Range("L36").Value = 35265
$ws.Range("N36").Value = -36389
$ws.Range("H46").Value = 3934.261
$ws.Range("I46").Value = 3499.25
$ws.Range("K46").Value = 3499.25
$ws.Range("M46").Value = -3311.25
$ws.Range("H132").Value = 5644.643
$ws.Range("I132").Value = 4596.2856
$ws.Range("K132").Value = 13788.8568
$ws.Range("M132").Value = -11258.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 79998
$ws.Range("J27").Value = 79998
$ws.Range("L27").Value = 79998
$ws.Range("N27").Value = -80136
$ws.Range("H115").Value = 79999
$ws.Range("J115").Value = 79999
$ws.Range("L115").Value = 79999
$ws.Range("N115").Value = -83133
$ws.Range("H126").Value = 850
$ws.Range("I126").Value = 700
$ws.Range("K126").Value = 2100
$ws.Range("M126").Value = 370
$ws.Range("H132").Value = 27024.684
$ws.Range("I132").Value = 2051.7334
$ws.Range("J132").Value = 95132.73
$ws.Range("K132").Value = 6155.2002
$ws.Range("L132").Value = 285398.19
$ws.Range("M132").Value = -3625.2002
$ws.Range("N132").Value = -290458.19
